$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 18:27"

# Albania and "Consejo Danes para los Refugiados" swapped position/order;
# update the two country-label cells to match the new pairing with their data rows
$ws.Range("A93").Value = "Albania"
$ws.Range("A94").Value = "Consejo Danes para los Refugiados"

# Updated COVID-19 statistics for the affected countries/rows
$ws.Range("B4").Value = 6440816
$ws.Range("C4").Value = 9664
$ws.Range("D4").Value = 3708171
$ws.Range("E4").Value = 2539689
$ws.Range("G4").Value = 138
$ws.Range("H4").Value = 192956
$ws.Range("B5").Value = 4196131
$ws.Range("C5").Value = 85292
$ws.Range("D5").Value = 3240977
$ws.Range("E5").Value = 883415
$ws.Range("G5").Value = 1060
$ws.Range("H5").Value = 71739
$ws.Range("E6").Value = 700032
$ws.Range("G6").Value = 36
$ws.Range("H6").Value = 126266
$ws.Range("B16").Value = 347152
$ws.Range("C16").Value = 2988
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 41551
$ws.Range("B21").Value = 279806
$ws.Range("C21").Value = 1578
$ws.Range("D21").Value = 251105
$ws.Range("E21").Value = 22028
$ws.Range("G21").Value = 53
$ws.Range("H21").Value = 6673
$ws.Range("B22").Value = 277634
$ws.Range("C22").Value = 1297
$ws.Range("D22").Value = 210015
$ws.Range("E22").Value = 32078
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 35541
$ws.Range("B24").Value = 251226
$ws.Range("C24").Value = 170
$ws.Range("E24").Value = 15617
$ws.Range("B72").Value = 27957
$ws.Range("C72").Value = 205
$ws.Range("D72").Value = 19049
$ws.Range("E72").Value = 8473
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 435
$ws.Range("B91").Value = 11524
$ws.Range("C91").Value = 138
$ws.Range("E91").Value = 7436
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 284
$ws.Range("B93").Value = 10255
$ws.Range("C93").Value = 153
$ws.Range("D93").Value = 6106
$ws.Range("E93").Value = 3833
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 316
$ws.Range("B94").Value = 10210
$ws.Range("C94").Value = 32
$ws.Range("D94").Value = 9439
$ws.Range("E94").Value = 511
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 260
$ws.Range("B124").Value = 3667
$ws.Range("C124").Value = 128
$ws.Range("E124").Value = 2062
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 41
$ws.Range("B140").Value = 2411
$ws.Range("C140").Value = 58
$ws.Range("D140").Value = 1756
$ws.Range("E140").Value = 639
$ws.Range("B146").Value = 2054
$ws.Range("C146").Value = 13
$ws.Range("D146").Value = 1611
$ws.Range("E146").Value = 372
$ws.Range("B181").Value = 337
$ws.Range("C181").Value = 1
$ws.Range("E181").Value = 1
